# Updates the cryptocurrency price/volume data to the latest scraped values
# (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.717.86"
$ws.Range("E2").Value = "  +1.57%  "

$ws.Range("D3").Value = "2.382.42"
$ws.Range("E3").Value = "  +3.78%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.56"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.35"
$ws.Range("E6").Value = "  +4.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.514"
$ws.Range("E7").Value = "  -4.45%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -0.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.27"
$ws.Range("E10").Value = "  +0.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.42"
$ws.Range("E11").Value = "  +2.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0815"
$ws.Range("E12").Value = "  -0.66%  "

$ws.Range("E13").Value = "  -0.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.03"
$ws.Range("E14").Value = "  -1.07%  "

$ws.Range("D15").Value = "2.752.10"
$ws.Range("E15").Value = "  +3.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.69"
$ws.Range("E16").Value = "  +4.51%  "

$ws.Range("D17").Value = "2.380.93"
$ws.Range("E17").Value = "  +4.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.816"
$ws.Range("E18").Value = "  +0.76%  "

$ws.Range("D19").Value = "43.630.70"
$ws.Range("E19").Value = "  +1.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.03"
$ws.Range("E20").Value = "  -3.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.32"
$ws.Range("E21").Value = "  +4.19%  "

$ws.Range("E22").Value = "  -0.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "241.92"
$ws.Range("E24").Value = "  +0.79%  "

$ws.Range("E25").Value = "  +2.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.63"
$ws.Range("E26").Value = "  +0.34%  "

$ws.Range("E27").Value = "  -0.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.81"
$ws.Range("E28").Value = "  +5.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.87"
$ws.Range("E29").Value = "  -2.67%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.20"
$ws.Range("E30").Value = "  -3.29%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.23"
$ws.Range("E31").Value = "  -4.20%  "

$ws.Range("E32").Value = "  -0.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "161.85"
$ws.Range("E33").Value = "  -3.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.29"
$ws.Range("E34").Value = "  -0.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.47"
$ws.Range("E35").Value = "  +4.36%  "

$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("E37").Value = "  +6.02%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.71"
$ws.Range("E38").Value = "  +11.76%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.10"
$ws.Range("E39").Value = "  -0.67%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0744"
$ws.Range("E40").Value = "  +0.71%  "

$ws.Range("E41").Value = "  +6.74%  "

$ws.Range("E42").Value = "  -0.96%  "

$ws.Range("E43").Value = "  -1.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.60"
$ws.Range("E44").Value = "  +13.51%  "

$ws.Range("D45").Value = "2.037.43"
$ws.Range("E45").Value = "  +3.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.75"
$ws.Range("E46").Value = "  +3.69%  "

$ws.Range("E47").Value = "  +0.78%  "

$ws.Range("E48").Value = "  +4.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.60"
$ws.Range("E49").Value = "  +7.77%  "

$ws.Range("E50").Value = "  +4.57%  "

$ws.Range("E51").Value = "  +0.92%  "

